$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.546.23'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.913.91'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.708'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +7.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '247.38'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.80'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.356'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.95'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0736'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0989'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.191.44'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.70'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.68%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.914.44'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.533.75'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.34'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.09%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.18'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '242.71'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.04%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.55'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.64'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.86'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.133'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.26'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0579'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.92'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +11.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.22'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.919'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.48'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +10.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.06'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +13.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.66'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.14'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.02%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0646'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.351.02'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.78%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.78'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.74'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.62%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.57'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Gas'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '12.17'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.88%  '
